$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22, shifting existing data rows 22-84 down to 23-85
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with a new daily price record
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44883
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112022
$ws.Range("G22").Value = "Arveja Verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 18000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 720
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
